# Simplify inputs and code, rely on prepDat()
#
# The English-language "translation" columns (man.source on Slurry,
# app.timing on Climate) are dropped - the workbook now keeps only the
# Danish values and the processing code (prepDat()) takes care of mapping
# them. The remaining header labels slide over to reuse the label that
# used to belong to the now-deleted column.

$wb = $excel.ActiveWorkbook

# --- Sheet "Slurry": drop column B (man.source: Pig/Cattle/Digestate) ---
$wsSlurry = $wb.Worksheets.Item("Slurry")
$wsSlurry.Range("B1").EntireColumn.Delete() | Out-Null
# Old column A header was "man.name"; reuse "man.source" for the remaining
# Danish manure-name column.
$wsSlurry.Range("A1").Value2 = "man.source"
$wsSlurry.Range("B9").Select() | Out-Null

# --- Sheet "Climate": drop column B (English month/season names) ---
$wsClimate = $wb.Worksheets.Item("Climate")
$wsClimate.Range("B1").EntireColumn.Delete() | Out-Null
# Old column A header was "app.timing.dk"; reuse "app.timing" for the
# remaining Danish timing column.
$wsClimate.Range("A1").Value2 = "app.timing"
$wsClimate.Range("C13").Select() | Out-Null

# --- Sheet "Application": no structural change, just keep selection/active ---
$wsApplication = $wb.Worksheets.Item("Application")
$wsApplication.Range("B7").Select() | Out-Null
$wsApplication.Activate() | Out-Null
